$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 791, shifting existing rows 791..832 down to 792..833
$ws.Rows.Item(791).Insert()

# Populate the newly inserted row 791 with the new data point.
# Column A holds a date-shaped string ("2026/02/14"); Excel's COM auto-detects
# such strings as real dates unless the cell is pre-formatted as Text. Force
# text formatting first, write the value, then reset the style to Normal so no
# residual number-format style lingers on the cell (matches its siblings).
$ws.Cells.Item(791, 1).NumberFormat = "@"
$ws.Cells.Item(791, 1).Value = "2026/02/14"
$ws.Cells.Item(791, 1).Style = "Normal"

$ws.Cells.Item(791, 2).Value = "土"
$ws.Cells.Item(791, 3).Value = 7
$ws.Cells.Item(791, 4).Value = 21
